$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 899
$ws.Range("I12").Value = 899
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 899
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -729
$ws.Range("N12").ClearContents()

$ws.Range("H18").Value = 732.8333
$ws.Range("I18").Value = 779.6
$ws.Range("J18").Value = 499
$ws.Range("K18").Value = 779.6
$ws.Range("L18").Value = 499
$ws.Range("M18").Value = -495.6
$ws.Range("N18").Value = -1067

$ws.Range("H21").Value = 20000
$ws.Range("J21").Value = 20000
$ws.Range("L21").Value = 20000
$ws.Range("N21").Value = -20936

$ws.Range("H23").Value = 20000
$ws.Range("J23").Value = 20000
$ws.Range("L23").Value = 20000
$ws.Range("N23").Value = -20468

$ws.Range("H32").Value = 1799.4117
$ws.Range("I32").Value = 1509.091
$ws.Range("K32").Value = 1509.091
$ws.Range("M32").Value = -1183.091

$ws.Range("H39").Value = 8
$ws.Range("I39").Value = 8
$ws.Range("K39").Value = 24
$ws.Range("M39").Value = 272

$ws.Range("H82").Value = 299.33334
$ws.Range("I82").Value = 299.33334
$ws.Range("K82").Value = 898.0000200000001
$ws.Range("M82").Value = -492.0000200000001

$ws.Range("H85").Value = 299.33334
$ws.Range("I85").Value = 299.33334
$ws.Range("K85").Value = 898.0000200000001
$ws.Range("M85").Value = 505.9999799999999

$ws.Range("H106").Value = 71459000
$ws.Range("I106").Value = 76950070
$ws.Range("K106").Value = 76950070
$ws.Range("M106").Value = -76949439

$ws.Range("H113").Value = 4588.3335
$ws.Range("I113").Value = 4487.5
$ws.Range("K113").Value = 4487.5
$ws.Range("M113").Value = -1233.5

$ws.Range("H116").Value = 11683
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 11683
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 11683
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -18567

$ws.Range("H137").Value = 1612.5714
$ws.Range("I137").Value = 1475.5294
$ws.Range("J137").Value = 2195
$ws.Range("K137").Value = 4426.5882
$ws.Range("L137").Value = 6585
$ws.Range("M137").Value = -1876.5882
$ws.Range("N137").Value = -11685

$ws.Range("H138").Value = 1813.2858
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6320.963
$ws.Range("I32").Value = 3245.2632
$ws.Range("J32").Value = 13625.75
$ws.Range("K32").Value = 3245.2632
$ws.Range("L32").Value = 13625.75
$ws.Range("M32").Value = -2958.2632
$ws.Range("N32").Value = -14199.75

$ws.Range("H44").Value = 10597.4
$ws.Range("J44").Value = 10597.4
$ws.Range("L44").Value = 10597.4
$ws.Range("N44").Value = -11573.4

$ws.Range("H55").Value = 19712.715
$ws.Range("J55").Value = 26388.2
$ws.Range("L55").Value = 26388.2
$ws.Range("N55").Value = -27018.2

$ws.Range("H61").Value = 1141.3334
$ws.Range("I61").Value = 969.6
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 969.6
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -757.6
$ws.Range("N61").Value = -2424

$ws.Range("H74").Value = 568
$ws.Range("I74").Value = 409.22223
$ws.Range("K74").Value = 409.22223
$ws.Range("M74").Value = 464.77777

$ws.Range("H77").Value = 568
$ws.Range("I77").Value = 409.22223
$ws.Range("K77").Value = 2046.11115
$ws.Range("M77").Value = 2321.88885

$ws.Range("H122").Value = 16663.545
$ws.Range("I122").Value = 8984.875
$ws.Range("J122").Value = 37140
$ws.Range("K122").Value = 26954.625
$ws.Range("L122").Value = 111420
$ws.Range("M122").Value = -24504.625
$ws.Range("N122").Value = -116320

$ws.Range("H132").Value = 1474.25
$ws.Range("I132").Value = 1399.1428
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 4197.428400000001
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -1667.428400000001
$ws.Range("N132").Value = -11060

$ws.Range("H135").Value = 200000
$ws.Range("J135").Value = 200000
$ws.Range("L135").Value = 200000
$ws.Range("N135").Value = -210140

$ws.Range("H136").Value = 1141.3334
$ws.Range("I136").Value = 969.6
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 2908.8
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -358.8000000000002
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 151.5
$ws.Range("I8").Value = 151.5
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 151.5
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -11.5
$ws.Range("N8").ClearContents()

$ws.Range("H80").Value = 281.77777
$ws.Range("I80").Value = 244.25
$ws.Range("K80").Value = 244.25
$ws.Range("M80").Value = 753.75

$ws.Range("H83").Value = 281.77777
$ws.Range("I83").Value = 244.25
$ws.Range("K83").Value = 1221.25
$ws.Range("M83").Value = 3770.75

$ws.Range("H107").Value = 1569.15
$ws.Range("I107").Value = 1234.9286
$ws.Range("J107").Value = 2349
$ws.Range("K107").Value = 1234.9286
$ws.Range("L107").Value = 2349
$ws.Range("M107").Value = 685.0714
$ws.Range("N107").Value = -6189

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 8714.5
$ws.Range("J15").Value = 3009
$ws.Range("L15").Value = 3009
$ws.Range("N15").Value = -3349

$ws.Range("H22").Value = 530
$ws.Range("I22").Value = 530
$ws.Range("K22").Value = 530
$ws.Range("M22").Value = -180

$ws.Range("H31").Value = 5127.4375
$ws.Range("I31").Value = 4078.7144
$ws.Range("K31").Value = 4078.7144
$ws.Range("M31").Value = -3783.7144

$ws.Range("H34").Value = 5127.4375
$ws.Range("I34").Value = 4078.7144
$ws.Range("K34").Value = 4078.7144
$ws.Range("M34").Value = -3876.7144

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 611.7714
$ws.Range("I4").Value = 565.76666
$ws.Range("K4").Value = 1697.29998
$ws.Range("M4").Value = -1585.29998

$ws.Range("H17").Value = 1431.25

$ws.Range("H62").Value = 2153.8462
$ws.Range("J62").Value = 2153.8462
$ws.Range("L62").Value = 6461.5386
$ws.Range("N62").Value = -7833.5386

$ws.Range("H65").Value = 2153.8462
$ws.Range("J65").Value = 2153.8462
$ws.Range("L65").Value = 19384.6158
$ws.Range("N65").Value = -26248.6158

$ws.Range("H69").Value = 2679.8
$ws.Range("I69").Value = 3500
$ws.Range("J69").Value = 2474.75
$ws.Range("K69").Value = 10500
$ws.Range("L69").Value = 7424.25
$ws.Range("M69").Value = -9689
$ws.Range("N69").Value = -9046.25

$ws.Range("H72").Value = 2679.8
$ws.Range("I72").Value = 3500
$ws.Range("J72").Value = 2474.75
$ws.Range("K72").Value = 31500
$ws.Range("L72").Value = 22272.75
$ws.Range("M72").Value = -27444
$ws.Range("N72").Value = -30384.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 49993.332
$ws.Range("J93").Value = 49993.332
$ws.Range("L93").Value = 49993.332
$ws.Range("N93").Value = -53737.332

$ws.Range("H122").Value = 43796.293
$ws.Range("I122").Value = 1805.85
$ws.Range("K122").Value = 5417.549999999999
$ws.Range("M122").Value = -2967.549999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 19152.5
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 19152.5
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 19152.5
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -19496.5

$ws.Range("H93").Value = 890.1053000000001
$ws.Range("I93").Value = 830
$ws.Range("J93").Value = 1401
$ws.Range("K93").Value = 830
$ws.Range("L93").Value = 1401
$ws.Range("M93").Value = 418
$ws.Range("N93").Value = -3897

$ws.Range("H122").Value = 1841.8572
$ws.Range("I122").Value = 1841.8572
$ws.Range("K122").Value = 5525.571599999999
$ws.Range("M122").Value = -3075.571599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 56666.668

$ws.Range("H73").Value = 56666.668
